# Update SnippetID values (column H) in the "Voice Lines - main" sheet.
# These are new, more stable SnippetIDs generated to replace the old ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$snippetIds = @{
    2  = "YkTdef"
    3  = "chW8xg"
    4  = "jZ2lVG"
    5  = "PXBE7j"
    6  = "DvVSQ6"
    7  = "qQOReF"
    8  = "qQOReF"
    9  = "qQOReF"
    10 = "qQOReF"
    11 = "qQOReF"
    12 = "iZ38Lw"
    13 = "32hKE1"
    14 = "U30VuF"
    15 = "i1GMUb"
    16 = "MdGHyj"
    17 = "MdGHyj"
    18 = "15ftGt"
    19 = "bbXZ4o"
    20 = "9YpLlu"
    21 = "9YpLlu"
    22 = "kRZMEF"
    23 = "AQWHIa"
    24 = "xo4coL"
    25 = "BZAGzy"
    26 = "IkrPq2"
    27 = "lTztHz"
    28 = "0pCnXe"
    29 = "GBmTaQ"
}

foreach ($row in $snippetIds.Keys) {
    $ws.Range("H$row").Value = $snippetIds[$row]
}
